$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new weekly data point sits between the existing rows that were
# D56=44252 (old row 56) and the prior row (old row 55, D=44488).
# Insert a new row at 56, pushing rows 56:131 down to 57:132.
$ws.Rows("56:56").Insert()

# Populate the newly inserted row 56 with the new "Orégano" price record.
$ws.Cells.Item(56,1).Value = 6
$ws.Cells.Item(56,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(56,3).Value = "Metropolitana"
$ws.Cells.Item(56,4).Value = 44546
$ws.Cells.Item(56,5).Value = 13
$ws.Cells.Item(56,6).Value = 100112029
$ws.Cells.Item(56,7).Value = "Orégano"
$ws.Cells.Item(56,8).Value = "Sin especificar"
$ws.Cells.Item(56,9).Value = "Primera"
$ws.Cells.Item(56,10).Value = 34
$ws.Cells.Item(56,11).Value = 8000
$ws.Cells.Item(56,12).Value = 9000
$ws.Cells.Item(56,13).Value = 8441
$ws.Cells.Item(56,14).Value = "`$/docena de atados"
$ws.Cells.Item(56,15).Value = "Región Metropolitana"
$ws.Cells.Item(56,16).Value = 2814
$ws.Cells.Item(56,17).Value = 3
$ws.Cells.Item(56,18).Value = "Hortaliza"
